$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update D176:D179 and C177:C180 (plain values, chained) ---
$ws.Range("D176").Value = 94.594999999999999

$ws.Range("C177").Value = 94.594999999999999
$ws.Range("D177").Value = 94.765000000000001

$ws.Range("C178").Value = 94.765000000000001
$ws.Range("D178").Value = 94.974999999999994

$ws.Range("C179").Value = 94.974999999999994
$ws.Range("D179").Value = 95.125

$ws.Range("C180").Value = 95.125
$ws.Range("D180").Value = 95.424999999999997

# --- D180 gets the new (non-bold, black) font style ---
$ws.Range("D180").Font.Color = 0
$ws.Range("D180").Font.Bold = $false

# --- C181 becomes a plain formula referencing D180, same new style ---
$ws.Range("C181").Formula = "=D180"
$ws.Range("C181").Font.Color = 0
$ws.Range("C181").Font.Bold = $false

$ws.Range("D181").Value = 95.564999999999998

# --- C182:C191 become a shared formula chain referencing the row above's D, same style ---
$ws.Range("C182:C191").Formula = "=D181"
$ws.Range("C182:C191").Font.Color = 0
$ws.Range("C182:C191").Font.Bold = $false

$ws.Range("D182").Value = 95.784999999999997
$ws.Range("D183").Value = 95.995000000000005
$ws.Range("D184").Value = 96.165000000000006
$ws.Range("D185").Value = 96.344999999999999
$ws.Range("D186").Value = 96.594999999999999
$ws.Range("D187").Value = 97.094999999999999
# D188:D191 values are unchanged (103.035, 103.295, 103.44499999999999, 103.62)

# --- View-state housekeeping (scroll position / zoom / selection) ---
$excel.ActiveWindow.ScrollRow = 133
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 111
$ws.Range("G153").Select() | Out-Null
